# Costos.20101018.xlsx - "Ajusté el tema costos."
# Bump the B4 growth formula by +2 (74 -> 76) and correct the C4/D4
# actuals for Sprint 2. The B5:B7 running totals are formula-driven
# (=prev+74) so they recompute automatically once B4 changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Formula = "=B3+76"
$ws.Range("C4").Value = 144
$ws.Range("D4").Value = 152

# Leave the active cell on D4, matching the editor's last touched cell.
$ws.Range("D4").Select()
